# Weekly data refresh: insert the newest week's record at the top of the
# data block (row 6, right after the 4 untouched historical rows 2-5) and
# push all the other weeks down by one row. The oldest row (previously 62)
# lands on the new row 63.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("6:6").Insert()

$ws.Cells.Item(6,1).Value = 6
$ws.Cells.Item(6,2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(6,3).Value = 'Metropolitana'
$ws.Cells.Item(6,4).Value = 44616
$ws.Cells.Item(6,5).Value = 13
$ws.Cells.Item(6,6).Value = 'Fruta'
$ws.Cells.Item(6,7).Value = 100108
$ws.Cells.Item(6,8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(6,9).Value = 100108007
$ws.Cells.Item(6,10).Value = 'Coco'
$ws.Cells.Item(6,11).Value = 'Sin especificar'
$ws.Cells.Item(6,12).Value = 'Primera'
$ws.Cells.Item(6,13).Value = 150
$ws.Cells.Item(6,14).Value = 22000
$ws.Cells.Item(6,15).Value = 22000
$ws.Cells.Item(6,16).Value = 22000
$ws.Cells.Item(6,17).Value = '$/malla 20 unidades'
$ws.Cells.Item(6,18).Value = 'Perú'
$ws.Cells.Item(6,19).Value = 1100
$ws.Cells.Item(6,20).Value = 20
